$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.286474823951721
$ws.Range("B1").Value = 2.763383150100708
$ws.Range("C1").Value = 4.658699512481689
$ws.Range("D1").Value = 2.023508071899414
$ws.Range("E1").Value = 1.263125658035278
